# Track liberos + set gender limits on each team
#
# Inserts a new "Libero" column before the existing "Setter" column (H),
# pushing Setter/Senior/Conflict-1..7 one column to the right, then marks
# the liberos (Holly Han - row 4, Ken Kirk - row 14) with "Y".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at H; everything from H onward shifts right by one.
$ws.Columns.Item(8).EntireColumn.Insert()

# New column header.
$ws.Range("H1").Value = "Libero"

# Mark the liberos.
$ws.Range("H4").Value = "Y"
$ws.Range("H14").Value = "Y"
